$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bug fix: Jung Ho Kang's position was wrong (2B) - should be SS ---
$ws.Cells.Item(169, 4).Value = "SS"

# --- Add status: new draft picks for "pk dodgers" ---
$ws.Cells.Item(318, 1).Value = "pk dodgers"
$ws.Cells.Item(318, 2).Value = "David Freese"
$ws.Cells.Item(318, 3).Value = 6
$ws.Cells.Item(318, 4).Value = "CI"

$ws.Cells.Item(319, 1).Value = "pk dodgers"
$ws.Cells.Item(319, 2).Value = "Derek Norris"
$ws.Cells.Item(319, 3).Value = 12
$ws.Cells.Item(319, 4).Value = "C"

$ws.Cells.Item(320, 1).Value = "pk dodgers"
$ws.Cells.Item(320, 2).Value = "Houston Street"
$ws.Cells.Item(320, 3).Value = 18
$ws.Cells.Item(320, 4).Value = "P"

# --- Update the view: scroll to where the edits were made and select D171 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 175
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D171").Select()
